$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.611
$ws.Range("D4").Value = -7.523000000000001
$ws.Range("D6").Value = -7.906000000000001
$ws.Range("A9").Value = -20.912
$ws.Range("D10").Value = -7.444
$ws.Range("B11").Value = 6.632000000000001
$ws.Range("D11").Value = -8.637
$ws.Range("A18").Value = -21.825
$ws.Range("A20").Value = -21.738
$ws.Range("C21").Value = -11.944
$ws.Range("D21").Value = -7.717000000000001
